# Append: 2025-10-11 12:32 JST
# Rebuild the data rows of the "ランサーズ" sheet with:
#  - all previous rows getting the new scrape timestamp
#  - three freshly scraped postings added in
#  - the whole list re-sorted by 優先度スコア (priority score, column G) descending
#  - widened D and H columns
#  - hyperlinks re-created on every F cell (F2:F10)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-10-11 12:32:06"

# Final row order/content (already sorted by column G descending).
$rows = @(
    @{ B = "【AI開発者募集】多機能転売ツールの構築をお願いします!"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5407785"; G = 420; H = "🔥AI,Ai ◆ツール,開発" },
    @{ B = "【急募】紙の伝票をWEBシステムへ自動データ入力開発【AI使用可能】"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411519"; G = 383; H = "🔥AI,Ai ◆開発" },
    @{ B = "急募 PR Zoom/Meet×TLDV×ChatGPT×Notion×Slack 議事録ワークフロー構築依頼"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5410688"; G = 323; H = "🔥GPT,ChatGPT" },
    @{ B = "【急募】仕事の予約システム構築をお手伝いください!【AI使用可能】"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411525"; G = 318; H = "🔥AI,Ai" },
    @{ B = "【急募】配送状況を自動取得するAPI開発者募集!"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411268"; G = 238; H = "🔥API ◆開発" },
    @{ B = "急募バックエンドエンジニア マッチングサイトの開発"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5404059"; G = 93; H = "◆開発 ◇サイト" },
    @{ B = "スプレッドシートをもとにした顧客・売上管理アプリのグライド化(Glide/無料版)"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411304"; G = 55; H = "◇アプリ" },
    @{ B = "【急募】時間単位で入札できるシステム構築の依頼"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411365"; G = 33; H = "" },
    @{ B = "【フォーム制作】物件見積り査定フォーム制作の依頼"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411435"; G = 13; H = "" }
)

# Remove any pre-existing hyperlinks before rewriting the range (avoids stale rId bindings).
$ws.Range("A1:H7").Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H -ne "") {
        $ws.Cells.Item($r, 8).Value = $row.H
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row.F) | Out-Null
    # Hyperlinks.Add stamps a fresh (but equivalent) style record; re-applying the
    # named style collapses it back onto the workbook's existing "Hyperlink" xf.
    $ws.Cells.Item($r, 6).Style = "Hyperlink"

    $r = $r + 1
}

# Column width tweaks (raw XML "width" runs ~0.83 narrower than the ColumnWidth property).
$ws.Columns.Item(4).ColumnWidth = 29.17
$ws.Columns.Item(8).ColumnWidth = 15.17
